$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert two new blank rows at position 8 (this pushes the existing
#    rows 8-12 down to 10-14: PANADOL, POWER B COMPLEX, the "sofa" row,
#    the totals row and the footer row).
# ------------------------------------------------------------------
$ws.Rows.Item(8).Insert(-4121, 0)
$ws.Rows.Item(8).Insert(-4121, 0)

# ------------------------------------------------------------------
# 2. Copy the formatting (styles/borders/number formats) of the
#    original data row (row 7) onto the two freshly inserted rows so
#    that they look identical to the other data rows.
# ------------------------------------------------------------------
$fmtSrc = $ws.Range("A7:Q7")
$fmtSrc.Copy()
$ws.Range("A8:Q9").PasteSpecial(-4122, 0, $false, $false)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3. Re-create the merged cells for the two new rows (A:B, C:G, H:K,
#    L:M, N:O) just like every other data row in the table.
# ------------------------------------------------------------------
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

# ------------------------------------------------------------------
# 4. Set the row heights to match the rest of the table.
# ------------------------------------------------------------------
$ws.Rows.Item(8).RowHeight = 24.75
$ws.Rows.Item(9).RowHeight = 25.5

# ------------------------------------------------------------------
# 5. Fill in the data for the two new rows. The "order limit" (L) and
#    "sale price" (P) columns use a numeric display format even though
#    the values are stored as plain text in this workbook, so the
#    number format is temporarily switched to Text while the value is
#    written and then restored to keep the original style.
# ------------------------------------------------------------------
function Set-TextValue($range, $value) {
    $orig = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = $orig
}

# Row 8 : GYPSUM  SYRUP
$ws.Cells.Item(8, 1).Value = 2
$ws.Range("C8").Value = "GYPSUM  SYRUP"
$ws.Range("H8").Value = "0:0"
Set-TextValue $ws.Range("L8") "0"
$ws.Range("N8").Value = "66.00"
Set-TextValue $ws.Range("P8") "66.0000"
$ws.Range("Q8").Value = "1:0"

# Row 9 : INDERAL 10 MG 50 TABS
$ws.Cells.Item(9, 1).Value = 3
$ws.Range("C9").Value = "INDERAL 10 MG 50 TABS"
$ws.Range("H9").Value = "0:1"
Set-TextValue $ws.Range("L9") "1"
$ws.Range("N9").Value = "75.00"
Set-TextValue $ws.Range("P9") "75.0000"
$ws.Range("Q9").Value = "1:0"

# ------------------------------------------------------------------
# 6. Renumber the "م" (index) column for the rows that shifted down.
# ------------------------------------------------------------------
$ws.Cells.Item(10, 1).Value = 4
$ws.Cells.Item(11, 1).Value = 5
$ws.Cells.Item(12, 1).Value = 6

# ------------------------------------------------------------------
# 7. Update the totals cell (now on row 13) to reflect the two new
#    rows that were added to the table.
# ------------------------------------------------------------------
$ws.Range("P13").Value = 251.31999999999999

# ------------------------------------------------------------------
# 8. Update the generation timestamp in the footer (now on row 14).
# ------------------------------------------------------------------
$ws.Range("A14").Value = "Thursday, 7 August, 2025 10:00 AM"
